# Generate Report for Handback
# Appends a new handback record (file 8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md)
# as row 4 on the "Overview", "zh-cn" and "de-de" worksheets, wires up the
# matching hyperlinks, and grows each sheet's table to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table "Overview" / A1:G3 -> A1:G4)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
$wsOverview.Range("B4").Value = "e2e\8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-26 14:56:38"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbaef77beb718cf5091ac36c2c52244116ca0906/e2e/8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md",
    $null,
    $null,
    "e2e\8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table "zh-cn" / A1:P3 -> A1:P4)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.cbaef77beb718cf5091ac36c2c52244116ca0906.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-26 14:56:33"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
$wsZhCn.Range("J4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.cbaef77beb718cf5091ac36c2c52244116ca0906.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-26 14:56:50"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbaef77beb718cf5091ac36c2c52244116ca0906/e2e/8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md",
    $null,
    $null,
    "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
) | Out-Null
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/cbaef77beb718cf5091ac36c2c52244116ca0906/e2e/8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md",
    $null,
    $null,
    "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" (table "de-de" / A1:P3 -> A1:P4)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.cbaef77beb718cf5091ac36c2c52244116ca0906.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-26 14:56:38"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
$wsDeDe.Range("J4").Value = "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.cbaef77beb718cf5091ac36c2c52244116ca0906.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-26 14:56:57"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbaef77beb718cf5091ac36c2c52244116ca0906/e2e/8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md",
    $null,
    $null,
    "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
) | Out-Null
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cbaef77beb718cf5091ac36c2c52244116ca0906/e2e/8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md",
    $null,
    $null,
    "8750d5d1-3cdf-4ffd-bf80-093f1a01723f.md"
) | Out-Null

# ---------------------------------------------------------------------
# Grow the three tables (ListObjects) so their range / autofilter / the
# worksheet dimension all cover the freshly written row 4.
# ---------------------------------------------------------------------
foreach ($ws in @($wsOverview, $wsZhCn, $wsDeDe)) {
    for ($i = 1; $i -le $ws.ListObjects.Count; $i++) {
        $lo = $ws.ListObjects.Item($i)
        if ($lo.Range.Rows.Count -lt 4) {
            $lo.Resize($lo.Range.Worksheet.Range($lo.Range.Cells.Item(1,1), $ws.Cells.Item(4, $lo.Range.Columns.Count)))
        }
    }
}
